# "class: adds new access code"
#
# Appends a "failure" run onto the existing "10/10: " entry and then logs
# three more dated entries: "10/12: application context", "10/17: " and
# "10/19: ".
#
# Plain Range.InsertAfter()/TypeText-style insertion silently folds new
# text into an adjacent run once it already shares that run's formatting,
# which would merge "10/10: " and "failure" (or "10/12: " and "application
# context") into one run instead of leaving them as sibling <w:r> elements.
# Splicing the literal WordprocessingML for each paragraph via
# Range.InsertXML keeps the runs distinct, matching how the entries were
# actually authored.

$d = $word.ActiveDocument
$wNs = "xmlns:w=`"http://schemas.openxmlformats.org/wordprocessingml/2006/main`""

# "10/10: " -> "10/10: " + "failure" (two runs in the existing paragraph)
$lastRange = $d.Paragraphs.Last.Range
$lastRange.InsertXML("<w:p $wNs>" +
    "<w:r><w:t xml:space=`"preserve`">10/10: </w:t></w:r>" +
    "<w:r><w:t>failure</w:t></w:r>" +
    "</w:p>")

# New paragraph: "10/12: " + "application context"
$p1 = $d.Paragraphs.Add()
$p1.Range.InsertXML("<w:p $wNs>" +
    "<w:r><w:t xml:space=`"preserve`">10/12: </w:t></w:r>" +
    "<w:r><w:t>application context</w:t></w:r>" +
    "</w:p>")

# New paragraph: "10/17: "
$p2 = $d.Paragraphs.Add()
$p2.Range.InsertXML("<w:p $wNs>" +
    "<w:r><w:t xml:space=`"preserve`">10/17: </w:t></w:r>" +
    "</w:p>")

# New paragraph: "10/19: "
$p3 = $d.Paragraphs.Add()
$p3.Range.InsertXML("<w:p $wNs>" +
    "<w:r><w:t xml:space=`"preserve`">10/19: </w:t></w:r>" +
    "</w:p>")
